# The deck currently has two themes:
#   ppt/theme/theme1.xml -> "Office Theme" (unused directly by the slide master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the one SlideMaster/Design)
#
# The target edit swaps the two themes' content, so the design actually
# applied to the slide master becomes the standard "Office Theme" colour
# scheme (and the previously-unused part becomes the "Integral" one).
#
# Re-colour the presentation's design (SlideMaster) theme to the Office
# Theme palette via the ThemeColorScheme, slot by slot. The twelve slots
# are, in order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

$themeColors.Item(1).RGB  = 0        # dk1      -> 000000
$themeColors.Item(2).RGB  = 16777215 # lt1      -> FFFFFF
$themeColors.Item(3).RGB  = 6968388  # dk2      -> 44546A
$themeColors.Item(4).RGB  = 15132391 # lt2      -> E7E6E6
$themeColors.Item(5).RGB  = 13998939 # accent1  -> 5B9BD5
$themeColors.Item(6).RGB  = 3243501  # accent2  -> ED7D31
$themeColors.Item(7).RGB  = 10855845 # accent3  -> A5A5A5
$themeColors.Item(8).RGB  = 49407    # accent4  -> FFC000
$themeColors.Item(9).RGB  = 12874308 # accent5  -> 4472C4
$themeColors.Item(10).RGB = 4697456  # accent6  -> 70AD47
$themeColors.Item(11).RGB = 12673797 # hlink    -> 0563C1
$themeColors.Item(12).RGB = 7491477  # folHlink -> 954F72
